# [MIG] 12.0 account_bank_statement_import_adyen, account_bank_statement_clearing_account
#
# The Adyen test workbook is updated so the sample statement is denominated
# in USD instead of EUR, and one of the sample "Gross Debit" amounts is
# corrected. The active selection on the sheet is also moved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Every cell whose text was "EUR" (Gross/Net Currency columns) becomes "USD".
# (Cells that already read e.g. "GBP" are left untouched.)
$ws.Cells.Replace("EUR", "USD")

# Row 10's Gross Debit (GC) amount changes from 666 to 1598.
$ws.Range("M10").Value = 1598

# Move the sheet's active selection.
$ws.Range("L9").Select()
